$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.70486980676651
$ws.Range("B1").Value = 3.114413499832153
$ws.Range("C1").Value = 2.863125801086426
$ws.Range("D1").Value = 2.368735551834106
$ws.Range("E1").Value = 2.109610319137573
